$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New token / count data for rows 2-18 (columns A, B, C)
$data = @(
    @("<people>",  "<people>", 26),
    @("<number>",  "<number>", 25),
    @("<time>",    "<come>",   16),
    @("<foxtrot>", "<foxtrot>",21),
    @("<write>",   "<mike>",   16),
    @("<been>",    "<been>",   13),
    @("<upward>",  "<upward>", 18),
    @("<escape>",  "<is>",     18),
    @("<he>",      "<see>",    23),
    @("<on>",      "<on>",     27),
    @("<there>",   "<there>",  21),
    @("<victor>",  "<six>",    18),
    @("<alt>",     "<alt>",    14),
    @("<as>",      "<as>",     19),
    @("<yankee>",  "<yankee>", 20),
    @("<him>",     "<six>",    18),
    @("<can>",     "<can>",    19)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
